# Edit script: apply betexplorer serie-c-group-b update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate match data (columns F:V) among shuffled row groups ---
# (Column A = index, column E = date stay fixed per physical row; only
#  the match-detail columns F..V were rotated among these rows.)

# Row 7
$ws.Range("F7").Value = 'Recanatese'
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 'Torres'
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 2.24
$ws.Range("K7").Value = '01/09/2023 09:12'
$ws.Range("L7").Value = 2.41
$ws.Range("M7").Value = '02/09/2023 11:26'
$ws.Range("N7").Value = 2.84
$ws.Range("O7").Value = '01/09/2023 09:12'
$ws.Range("P7").Value = 2.99
$ws.Range("Q7").Value = '02/09/2023 18:47'
$ws.Range("R7").Value = 3.3
$ws.Range("S7").Value = '01/09/2023 09:12'
$ws.Range("T7").Value = 3.21
$ws.Range("U7").Value = '02/09/2023 11:26'
$ws.Range("V7").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/recanatese-sassari-torres/6PNcuOH8/'

# Row 8
$ws.Range("F8").Value = 'Pescara'
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 'Juventus U23'
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1.74
$ws.Range("K8").Value = '01/09/2023 18:13'
$ws.Range("L8").Value = 2.03
$ws.Range("M8").Value = '02/09/2023 20:35'
$ws.Range("N8").Value = 3.66
$ws.Range("O8").Value = '01/09/2023 18:13'
$ws.Range("P8").Value = 3.34
$ws.Range("Q8").Value = '02/09/2023 20:35'
$ws.Range("R8").Value = 3.92
$ws.Range("S8").Value = '01/09/2023 18:13'
$ws.Range("T8").Value = 3.77
$ws.Range("U8").Value = '02/09/2023 20:25'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/pescara-juventus/vJVBHKP1/'

# Row 9
$ws.Range("F9").Value = 'Entella'
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 'Ancona'
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1.62
$ws.Range("K9").Value = '01/09/2023 09:12'
$ws.Range("L9").Value = 1.7
$ws.Range("M9").Value = '02/09/2023 16:50'
$ws.Range("N9").Value = 3.55
$ws.Range("O9").Value = '01/09/2023 09:12'
$ws.Range("P9").Value = 3.57
$ws.Range("Q9").Value = '02/09/2023 18:49'
$ws.Range("R9").Value = 4.9
$ws.Range("S9").Value = '01/09/2023 09:12'
$ws.Range("T9").Value = 4.89
$ws.Range("U9").Value = '02/09/2023 16:50'
$ws.Range("V9").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/entella-ancona/z5BAxM1R/'

# Row 10
$ws.Range("F10").Value = 'Spal'
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 'Vis Pesaro'
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1.34
$ws.Range("K10").Value = '01/09/2023 18:13'
$ws.Range("L10").Value = 1.35
$ws.Range("M10").Value = '02/09/2023 10:46'
$ws.Range("N10").Value = 4.23
$ws.Range("O10").Value = '01/09/2023 18:13'
$ws.Range("P10").Value = 4.67
$ws.Range("Q10").Value = '02/09/2023 20:39'
$ws.Range("R10").Value = 8.35
$ws.Range("S10").Value = '01/09/2023 18:13'
$ws.Range("T10").Value = 9.98
$ws.Range("U10").Value = '02/09/2023 20:39'
$ws.Range("V10").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/spal-vis-pesaro/SvC6w2nL/'

# Row 67
$ws.Range("F67").Value = 'Arezzo'
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 'Cesena'
$ws.Range("I67").Value = 2
$ws.Range("J67").Value = 4.99
$ws.Range("K67").Value = '05/10/2023 15:12'
$ws.Range("L67").Value = 6.14
$ws.Range("M67").Value = '09/10/2023 20:43'
$ws.Range("N67").Value = 3.49
$ws.Range("O67").Value = '05/10/2023 15:12'
$ws.Range("P67").Value = 4.21
$ws.Range("Q67").Value = '09/10/2023 20:43'
$ws.Range("R67").Value = 1.65
$ws.Range("S67").Value = '05/10/2023 15:12'
$ws.Range("T67").Value = 1.52
$ws.Range("U67").Value = '09/10/2023 20:43'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/arezzo-cesena/AaNV4ysk/'

# Row 69
$ws.Range("F69").Value = 'Gubbio'
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 'Carrarese'
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2.49
$ws.Range("K69").Value = '05/10/2023 15:12'
$ws.Range("L69").Value = 2.32
$ws.Range("M69").Value = '09/10/2023 19:22'
$ws.Range("N69").Value = 2.83
$ws.Range("O69").Value = '05/10/2023 15:12'
$ws.Range("P69").Value = 2.9
$ws.Range("Q69").Value = '09/10/2023 20:31'
$ws.Range("R69").Value = 2.88
$ws.Range("S69").Value = '05/10/2023 15:12'
$ws.Range("T69").Value = 3.54
$ws.Range("U69").Value = '09/10/2023 19:22'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/gubbio-carrarese/U1wX2FB1/'

# Row 78
$ws.Range("F78").Value = 'Pescara'
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 'Vis Pesaro'
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.48
$ws.Range("K78").Value = '12/10/2023 08:13'
$ws.Range("L78").Value = 1.29
$ws.Range("M78").Value = '16/10/2023 19:56'
$ws.Range("N78").Value = 4.03
$ws.Range("O78").Value = '12/10/2023 08:13'
$ws.Range("P78").Value = 5.11
$ws.Range("Q78").Value = '16/10/2023 20:43'
$ws.Range("R78").Value = 5.63
$ws.Range("S78").Value = '12/10/2023 08:13'
$ws.Range("T78").Value = 9.14
$ws.Range("U78").Value = '16/10/2023 20:43'
$ws.Range("V78").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/pescara-vis-pesaro/Iya0wdmE/'

# Row 79
$ws.Range("F79").Value = 'Pontedera'
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 'Rimini'
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 2.01
$ws.Range("K79").Value = '12/10/2023 08:13'
$ws.Range("L79").Value = 2.09
$ws.Range("M79").Value = '16/10/2023 20:35'
$ws.Range("N79").Value = 3.11
$ws.Range("O79").Value = '12/10/2023 08:13'
$ws.Range("P79").Value = 3.12
$ws.Range("Q79").Value = '16/10/2023 20:35'
$ws.Range("R79").Value = 3.54
$ws.Range("S79").Value = '12/10/2023 08:13'
$ws.Range("T79").Value = 3.88
$ws.Range("U79").Value = '16/10/2023 20:41'
$ws.Range("V79").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/us-pontedera-rimini/dSc4xG2K/'

# Row 80
$ws.Range("F80").Value = 'Torres'
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 'Pontedera'
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1.84
$ws.Range("K80").Value = '19/10/2023 08:13'
$ws.Range("L80").Value = 1.99
$ws.Range("M80").Value = '22/10/2023 13:53'
$ws.Range("N80").Value = 3.18
$ws.Range("O80").Value = '19/10/2023 08:13'
$ws.Range("P80").Value = 3.12
$ws.Range("Q80").Value = '22/10/2023 13:53'
$ws.Range("R80").Value = 4.12
$ws.Range("S80").Value = '19/10/2023 08:13'
$ws.Range("T80").Value = 4.29
$ws.Range("U80").Value = '22/10/2023 13:53'
$ws.Range("V80").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/sassari-torres-us-pontedera/rmbIqX9m/'

# Row 82
$ws.Range("F82").Value = 'Fermana'
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 'Entella'
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = 5.14
$ws.Range("K82").Value = '19/10/2023 08:13'
$ws.Range("L82").Value = 5
$ws.Range("M82").Value = '22/10/2023 13:52'
$ws.Range("N82").Value = 3.29
$ws.Range("O82").Value = '19/10/2023 08:13'
$ws.Range("P82").Value = 3.27
$ws.Range("Q82").Value = '22/10/2023 13:52'
$ws.Range("R82").Value = 1.68
$ws.Range("S82").Value = '19/10/2023 08:13'
$ws.Range("T82").Value = 1.81
$ws.Range("U82").Value = '22/10/2023 13:52'
$ws.Range("V82").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-entella/djxmUiA0/'

# Row 84
$ws.Range("F84").Value = 'Lucchese'
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 'Pescara'
$ws.Range("I84").Value = 4
$ws.Range("J84").Value = 2.5
$ws.Range("K84").Value = '19/10/2023 08:13'
$ws.Range("L84").Value = 2.62
$ws.Range("M84").Value = '22/10/2023 18:21'
$ws.Range("N84").Value = 3.35
$ws.Range("O84").Value = '19/10/2023 08:13'
$ws.Range("P84").Value = 3.56
$ws.Range("Q84").Value = '22/10/2023 18:24'
$ws.Range("R84").Value = 2.55
$ws.Range("S84").Value = '19/10/2023 08:13'
$ws.Range("T84").Value = 2.53
$ws.Range("U84").Value = '22/10/2023 18:21'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/lucchese-pescara/W8UaRkfJ/'

# Row 85
$ws.Range("F85").Value = 'Pineto'
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 'Cesena'
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = 5.08
$ws.Range("K85").Value = '19/10/2023 08:13'
$ws.Range("L85").Value = 5.06
$ws.Range("M85").Value = '22/10/2023 18:26'
$ws.Range("N85").Value = 3.54
$ws.Range("O85").Value = '19/10/2023 08:13'
$ws.Range("P85").Value = 3.53
$ws.Range("Q85").Value = '22/10/2023 18:26'
$ws.Range("R85").Value = 1.63
$ws.Range("S85").Value = '19/10/2023 08:13'
$ws.Range("T85").Value = 1.74
$ws.Range("U85").Value = '22/10/2023 18:26'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/pineto-cesena/pCIsjZWJ/'

# Row 86
$ws.Range("F86").Value = 'Rimini'
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 'Ancona'
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = 2.18
$ws.Range("K86").Value = '19/10/2023 08:13'
$ws.Range("L86").Value = 2.13
$ws.Range("M86").Value = '22/10/2023 18:21'
$ws.Range("N86").Value = 3.04
$ws.Range("O86").Value = '19/10/2023 08:13'
$ws.Range("P86").Value = 3.61
$ws.Range("Q86").Value = '22/10/2023 18:21'
$ws.Range("R86").Value = 3.18
$ws.Range("S86").Value = '19/10/2023 08:13'
$ws.Range("T86").Value = 3.21
$ws.Range("U86").Value = '22/10/2023 18:21'
$ws.Range("V86").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/rimini-ancona/O2HokgnQ/'

# Row 87
$ws.Range("F87").Value = 'Vis Pesaro'
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 'Recanatese'
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 2.81
$ws.Range("K87").Value = '19/10/2023 08:13'
$ws.Range("L87").Value = 2.59
$ws.Range("M87").Value = '22/10/2023 18:21'
$ws.Range("N87").Value = 2.9
$ws.Range("O87").Value = '19/10/2023 08:13'
$ws.Range("P87").Value = 2.98
$ws.Range("Q87").Value = '22/10/2023 18:21'
$ws.Range("R87").Value = 2.5
$ws.Range("S87").Value = '19/10/2023 08:13'
$ws.Range("T87").Value = 2.97
$ws.Range("U87").Value = '22/10/2023 18:21'
$ws.Range("V87").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/vis-pesaro-recanatese/xAcMriPg/'

# Row 95
$ws.Range("F95").Value = 'Spal'
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 'Sestri Levante'
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 1.76
$ws.Range("K95").Value = '24/10/2023 21:12'
$ws.Range("L95").Value = 1.85
$ws.Range("M95").Value = '26/10/2023 18:29'
$ws.Range("N95").Value = 3.21
$ws.Range("O95").Value = '24/10/2023 21:12'
$ws.Range("P95").Value = 3.03
$ws.Range("Q95").Value = '26/10/2023 18:29'
$ws.Range("R95").Value = 4.72
$ws.Range("S95").Value = '24/10/2023 21:12'
$ws.Range("T95").Value = 5.29
$ws.Range("U95").Value = '26/10/2023 18:29'
$ws.Range("V95").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/spal-sestri-levante/WW8ctJKh/'

# Row 96
$ws.Range("F96").Value = 'Pontedera'
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 'Vis Pesaro'
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1.84
$ws.Range("K96").Value = '24/10/2023 21:12'
$ws.Range("L96").Value = 1.82
$ws.Range("M96").Value = '26/10/2023 18:02'
$ws.Range("N96").Value = 3.08
$ws.Range("O96").Value = '24/10/2023 21:12'
$ws.Range("P96").Value = 3.34
$ws.Range("Q96").Value = '26/10/2023 18:02'
$ws.Range("R96").Value = 4.28
$ws.Range("S96").Value = '24/10/2023 21:12'
$ws.Range("T96").Value = 4.8
$ws.Range("U96").Value = '26/10/2023 18:02'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/us-pontedera-vis-pesaro/McJNnLsP/'

# Row 97
$ws.Range("F97").Value = 'Rimini'
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 'Lucchese'
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3.08
$ws.Range("K97").Value = '24/10/2023 21:12'
$ws.Range("L97").Value = 2.92
$ws.Range("M97").Value = '26/10/2023 18:26'
$ws.Range("N97").Value = 2.84
$ws.Range("O97").Value = '24/10/2023 21:12'
$ws.Range("P97").Value = 3.2
$ws.Range("Q97").Value = '26/10/2023 18:26'
$ws.Range("R97").Value = 2.41
$ws.Range("S97").Value = '24/10/2023 21:12'
$ws.Range("T97").Value = 2.48
$ws.Range("U97").Value = '26/10/2023 18:26'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/rimini-lucchese/htAgsa5n/'

# --- Step 2: append 9 new match rows (98..106 -> sheet rows 99..107) ---
# Clone formatting (styles) of the last existing row down across the new rows first
$ws.Range("A98:V98").Copy($ws.Range("A99:V107"))

# Sheet row 99 (Indice 98)
$ws.Range("A99").Value = 98
$ws.Range("B99").Value = 'italy'
$ws.Range("C99").Value = 'serie-c-group-b'
$ws.Range("D99").Value = '2023-2024'
$ws.Range("E99").Value = 45228.58333333334
$ws.Range("F99").Value = 'Juventus U23'
$ws.Range("G99").Value = 3
$ws.Range("H99").Value = 'Olbia'
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 2.08
$ws.Range("K99").Value = '26/10/2023 22:12'
$ws.Range("L99").Value = 2.06
$ws.Range("M99").Value = '29/10/2023 13:33'
$ws.Range("N99").Value = 2.94
$ws.Range("O99").Value = '26/10/2023 22:12'
$ws.Range("P99").Value = 3.17
$ws.Range("Q99").Value = '29/10/2023 13:50'
$ws.Range("R99").Value = 3.57
$ws.Range("S99").Value = '26/10/2023 22:12'
$ws.Range("T99").Value = 3.6
$ws.Range("U99").Value = '29/10/2023 13:33'
$ws.Range("V99").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/juventus-olbia/Q7ezL8lp/'

# Sheet row 100 (Indice 99)
$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 'italy'
$ws.Range("C100").Value = 'serie-c-group-b'
$ws.Range("D100").Value = '2023-2024'
$ws.Range("E100").Value = 45228.58333333334
$ws.Range("F100").Value = 'Torres'
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 'Spal'
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 1.79
$ws.Range("K100").Value = '26/10/2023 22:12'
$ws.Range("L100").Value = 2
$ws.Range("M100").Value = '29/10/2023 13:53'
$ws.Range("N100").Value = 3.12
$ws.Range("O100").Value = '26/10/2023 22:12'
$ws.Range("P100").Value = 3.18
$ws.Range("Q100").Value = '29/10/2023 13:53'
$ws.Range("R100").Value = 4.5
$ws.Range("S100").Value = '26/10/2023 22:12'
$ws.Range("T100").Value = 4.12
$ws.Range("U100").Value = '29/10/2023 13:54'
$ws.Range("V100").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/sassari-torres-spal/Q3uZuAQ9/'

# Sheet row 101 (Indice 100)
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 'italy'
$ws.Range("C101").Value = 'serie-c-group-b'
$ws.Range("D101").Value = '2023-2024'
$ws.Range("E101").Value = 45228.77083333334
$ws.Range("F101").Value = 'Pescara'
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 'Recanatese'
$ws.Range("I101").Value = 3
$ws.Range("J101").Value = 1.59
$ws.Range("K101").Value = '26/10/2023 22:12'
$ws.Range("L101").Value = 1.44
$ws.Range("M101").Value = '29/10/2023 11:26'
$ws.Range("N101").Value = 3.9
$ws.Range("O101").Value = '26/10/2023 22:12'
$ws.Range("P101").Value = 4.52
$ws.Range("Q101").Value = '29/10/2023 18:19'
$ws.Range("R101").Value = 4.63
$ws.Range("S101").Value = '26/10/2023 22:12'
$ws.Range("T101").Value = 6.2
$ws.Range("U101").Value = '29/10/2023 16:58'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/pescara-recanatese/W4nQsWec/'

# Sheet row 102 (Indice 101)
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 'italy'
$ws.Range("C102").Value = 'serie-c-group-b'
$ws.Range("D102").Value = '2023-2024'
$ws.Range("E102").Value = 45228.86458333334
$ws.Range("F102").Value = 'Fermana'
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 'Ancona'
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = 2.79
$ws.Range("K102").Value = '26/10/2023 22:12'
$ws.Range("L102").Value = 3.54
$ws.Range("M102").Value = '29/10/2023 20:36'
$ws.Range("N102").Value = 2.8
$ws.Range("O102").Value = '26/10/2023 22:12'
$ws.Range("P102").Value = 3.01
$ws.Range("Q102").Value = '29/10/2023 20:36'
$ws.Range("R102").Value = 2.6
$ws.Range("S102").Value = '26/10/2023 22:12'
$ws.Range("T102").Value = 2.26
$ws.Range("U102").Value = '29/10/2023 20:36'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-ancona/pAs6xyKH/'

# Sheet row 103 (Indice 102)
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = 'italy'
$ws.Range("C103").Value = 'serie-c-group-b'
$ws.Range("D103").Value = '2023-2024'
$ws.Range("E103").Value = 45228.86458333334
$ws.Range("F103").Value = 'Lucchese'
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 'Pontedera'
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 2.01
$ws.Range("K103").Value = '26/10/2023 22:12'
$ws.Range("L103").Value = 1.93
$ws.Range("M103").Value = '29/10/2023 20:41'
$ws.Range("N103").Value = 3.09
$ws.Range("O103").Value = '26/10/2023 22:12'
$ws.Range("P103").Value = 3.59
$ws.Range("Q103").Value = '29/10/2023 20:42'
$ws.Range("R103").Value = 3.56
$ws.Range("S103").Value = '26/10/2023 22:12'
$ws.Range("T103").Value = 3.86
$ws.Range("U103").Value = '29/10/2023 20:42'
$ws.Range("V103").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/lucchese-us-pontedera/4bfvKS3j/'

# Sheet row 104 (Indice 103)
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 'italy'
$ws.Range("C104").Value = 'serie-c-group-b'
$ws.Range("D104").Value = '2023-2024'
$ws.Range("E104").Value = 45229.86458333334
$ws.Range("F104").Value = 'Arezzo'
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 'Gubbio'
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2.34
$ws.Range("K104").Value = '26/10/2023 22:12'
$ws.Range("L104").Value = 3.4
$ws.Range("M104").Value = '30/10/2023 20:42'
$ws.Range("N104").Value = 2.85
$ws.Range("O104").Value = '26/10/2023 22:12'
$ws.Range("P104").Value = 2.96
$ws.Range("Q104").Value = '30/10/2023 20:42'
$ws.Range("R104").Value = 3.09
$ws.Range("S104").Value = '26/10/2023 22:12'
$ws.Range("T104").Value = 2.35
$ws.Range("U104").Value = '30/10/2023 20:42'
$ws.Range("V104").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/arezzo-gubbio/xpvbvck5/'

# Sheet row 105 (Indice 104)
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = 'italy'
$ws.Range("C105").Value = 'serie-c-group-b'
$ws.Range("D105").Value = '2023-2024'
$ws.Range("E105").Value = 45229.86458333334
$ws.Range("F105").Value = 'Cesena'
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 'Carrarese'
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.84
$ws.Range("K105").Value = '26/10/2023 21:12'
$ws.Range("L105").Value = 1.65
$ws.Range("M105").Value = '30/10/2023 20:42'
$ws.Range("N105").Value = 3.08
$ws.Range("O105").Value = '26/10/2023 21:12'
$ws.Range("P105").Value = 3.44
$ws.Range("Q105").Value = '30/10/2023 20:42'
$ws.Range("R105").Value = 4.28
$ws.Range("S105").Value = '26/10/2023 21:12'
$ws.Range("T105").Value = 6.35
$ws.Range("U105").Value = '30/10/2023 20:42'
$ws.Range("V105").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/cesena-carrarese/W6w2wH4B/'

# Sheet row 106 (Indice 105)
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = 'italy'
$ws.Range("C106").Value = 'serie-c-group-b'
$ws.Range("D106").Value = '2023-2024'
$ws.Range("E106").Value = 45229.86458333334
$ws.Range("F106").Value = 'Perugia'
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 'Entella'
$ws.Range("I106").Value = 1
$ws.Range("J106").Value = 2.12
$ws.Range("K106").Value = '26/10/2023 22:12'
$ws.Range("L106").Value = 1.92
$ws.Range("M106").Value = '30/10/2023 20:36'
$ws.Range("N106").Value = 3.01
$ws.Range("O106").Value = '26/10/2023 22:12'
$ws.Range("P106").Value = 3.28
$ws.Range("Q106").Value = '30/10/2023 20:36'
$ws.Range("R106").Value = 3.36
$ws.Range("S106").Value = '26/10/2023 22:12'
$ws.Range("T106").Value = 4.34
$ws.Range("U106").Value = '30/10/2023 20:36'
$ws.Range("V106").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/perugia-entella/hEmMrCti/'

# Sheet row 107 (Indice 106)
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = 'italy'
$ws.Range("C107").Value = 'serie-c-group-b'
$ws.Range("D107").Value = '2023-2024'
$ws.Range("E107").Value = 45229.86458333334
$ws.Range("F107").Value = 'Vis Pesaro'
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 'Pineto'
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.21
$ws.Range("K107").Value = '26/10/2023 22:12'
$ws.Range("L107").Value = 2.65
$ws.Range("M107").Value = '30/10/2023 20:36'
$ws.Range("N107").Value = 2.88
$ws.Range("O107").Value = '26/10/2023 22:12'
$ws.Range("P107").Value = 2.76
$ws.Range("Q107").Value = '30/10/2023 20:36'
$ws.Range("R107").Value = 3.31
$ws.Range("S107").Value = '26/10/2023 22:12'
$ws.Range("T107").Value = 3.16
$ws.Range("U107").Value = '30/10/2023 20:36'
$ws.Range("V107").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-b/vis-pesaro-pineto/KQvwuUuG/'

Write-Host "Edit complete"